# Apply CountryConcordance.xlsx edits: add third "Exogenous country name" column
# with AQUASTAT <-> IFs <-> Exogenous country-name concordance data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    ,@('AQUASTAT country name','Country Name in IFs','Exogenous country name')
    ,@('Afghanistan','Afghanistan','Afghanistan')
    ,@('Albania','Albania','Albania')
    ,@('Algeria','Algeria','Algeria')
    ,@('Andorra','Not in IFs','#N/A')
    ,@('Angola','Angola','Angola')
    ,@('Antigua and Barbuda','Not in IFs','#N/A')
    ,@('Argentina','Argentina','Argentina')
    ,@('Armenia','Armenia','Armenia')
    ,@('Australia','Australia','Australia')
    ,@('Austria','Austria','Austria')
    ,@('Azerbaijan','Azerbaijan','Azerbaijan')
    ,@('Bahamas','Bahamas, The','Bahamas')
    ,@('Bahrain','Bahrain','Bahrain')
    ,@('Bangladesh','Bangladesh','Bangladesh')
    ,@('Barbados','Barbados','Barbados')
    ,@('Belarus','Belarus','Belarus')
    ,@('Belgium','Belgium','Belgium')
    ,@('Belize','Belize','Belize')
    ,@('Benin','Benin','Benin')
    ,@('Bhutan','Bhutan','Bhutan')
    ,@('Bolivia (Plurinational State of)','Bolivia','Bolivia')
    ,@('Bosnia and Herzegovina','Bosnia and Herzegovina','Bosnia')
    ,@('Botswana','Botswana','Botswana')
    ,@('Brazil','Brazil','Brazil')
    ,@('Brunei Darussalam','Brunei','Brunei')
    ,@('Bulgaria','Bulgaria','Bulgaria')
    ,@('Burkina Faso','Burkina Faso','Burkina Faso')
    ,@('Burundi','Burundi','Burundi')
    ,@('Cabo Verde','Cape Verde','Cape Verde')
    ,@('Cambodia','Cambodia','Cambodia')
    ,@('Cameroon','Cameroon','Cameroon')
    ,@('Canada','Canada','Canada')
    ,@('Central African Republic','Central African Republic','Central AfR')
    ,@('Chad','Chad','Chad')
    ,@('Chile','Chile','Chile')
    ,@('China','China','China')
    ,@('Colombia','Colombia','Colombia')
    ,@('Comoros','Comoros','Comoros')
    ,@('Congo','Congo, Republic of','Congo; Republic of')
    ,@('Cook Islands','Not in IFs','#N/A')
    ,@('Costa Rica','Costa Rica','Costa Rica')
    ,@('Côte d''Ivoire','Cote d''Ivoire','Cote d''Ivoire')
    ,@('Croatia','Croatia','Croatia')
    ,@('Cuba','Cuba','Cuba')
    ,@('Cyprus','Cyprus','Cyprus')
    ,@('Czechia','Czech Republic','Czech Republic')
    ,@('Democratic People''s Republic of Korea','Korea, Democratic People''s Republic of','Korea North')
    ,@('Democratic Republic of the Congo','Congo, Democratic Republic of','Congo; Democratic Republic of')
    ,@('Denmark','Denmark','Denmark')
    ,@('Djibouti','Djibouti','Djibouti')
    ,@('Dominica','Not in IFs','#N/A')
    ,@('Dominican Republic','Dominican Republic','DominicanRep')
    ,@('Ecuador','Ecuador','Ecuador')
    ,@('Egypt','Egypt, Arab Republic of','Egypt')
    ,@('El Salvador','El Salvador','El Salvador')
    ,@('Equatorial Guinea','Equatorial Guinea','Equa Guinea')
    ,@('Eritrea','Eritrea','Eritrea')
    ,@('Estonia','Estonia','Estonia')
    ,@('Ethiopia','Ethiopia','Ethiopia')
    ,@('Faroe Islands','Not in IFs','#N/A')
    ,@('Fiji','Fiji','Fiji')
    ,@('Finland','Finland','Finland')
    ,@('France','France','France')
    ,@('Gabon','Gabon','Gabon')
    ,@('Gambia','Gambia, The','Gambia')
    ,@('Georgia','Georgia','Georgia')
    ,@('Germany','Germany','Germany')
    ,@('Ghana','Ghana','Ghana')
    ,@('Greece','Greece','Greece')
    ,@('Grenada','Grenada','Grenada')
    ,@('Guatemala','Guatemala','Guatemala')
    ,@('Guinea','Guinea','Guinea')
    ,@('Guinea-Bissau','Guinea-Bissau','GuineaBiss')
    ,@('Guyana','Guyana','Guyana')
    ,@('Haiti','Haiti','Haiti')
    ,@('Holy See','Not in IFs','#N/A')
    ,@('Honduras','Honduras','Honduras')
    ,@('Hungary','Hungary','Hungary')
    ,@('Iceland','Iceland','Iceland')
    ,@('India','India','India')
    ,@('Indonesia','Indonesia','Indonesia')
    ,@('Iran (Islamic Republic of)','Iran, Islamic Republic of','Iran')
    ,@('Iraq','Iraq','Iraq')
    ,@('Ireland','Ireland','Ireland')
    ,@('Israel','Israel','Israel')
    ,@('Italy','Italy','Italy')
    ,@('Jamaica','Jamaica','Jamaica')
    ,@('Japan','Japan','Japan')
    ,@('Jordan','Jordan','Jordan')
    ,@('Kazakhstan','Kazakhstan','Kazakhstan')
    ,@('Kenya','Kenya','Kenya')
    ,@('Kiribati','Not in IFs','#N/A')
    ,@('Kuwait','Kuwait','Kuwait')
    ,@('Kyrgyzstan','Kyrgyz Republic','Kyrgyz')
    ,@('Lao People''s Democratic Republic','Laos, People''s Democratic Republic','Laos')
    ,@('Latvia','Latvia','Latvia')
    ,@('Lebanon','Lebanon','Lebanon')
    ,@('Lesotho','Lesotho','Lesotho')
    ,@('Liberia','Liberia','Liberia')
    ,@('Libya','Libya','Libya')
    ,@('Liechtenstein','Not in IFs','#N/A')
    ,@('Lithuania','Lithuania','Lithuania')
    ,@('Luxembourg','Luxembourg','Luxembourg')
    ,@('Madagascar','Madagascar','Madagascar')
    ,@('Malawi','Malawi','Malawi')
    ,@('Malaysia','Malaysia','Malaysia')
    ,@('Maldives','Maldives','Maldives')
    ,@('Mali','Mali','Mali')
    ,@('Malta','Malta','Malta')
    ,@('Marshall Islands','Not in IFs','#N/A')
    ,@('Mauritania','Mauritania','Mauritania')
    ,@('Mauritius','Mauritius','Mauritius')
    ,@('Mexico','Mexico','Mexico')
    ,@('Micronesia (Federated States of)','Micronesia (Federated States of)','Micronesia; Fed. Sts.')
    ,@('Monaco','Not in IFs','#N/A')
    ,@('Mongolia','Mongolia','Mongolia')
    ,@('Montenegro','Montenegro','Montenegro')
    ,@('Morocco','Morocco','Morocco')
    ,@('Mozambique','Mozambique','Mozambique')
    ,@('Myanmar','Myanmar','Myanmar')
    ,@('Namibia','Namibia','Namibia')
    ,@('Nauru','Not in IFs','#N/A')
    ,@('Nepal','Nepal','Nepal')
    ,@('Netherlands','Netherlands','Netherlands')
    ,@('New Zealand','New Zealand','New Zealand')
    ,@('Nicaragua','Nicaragua','Nicaragua')
    ,@('Niger','Niger','Niger')
    ,@('Nigeria','Nigeria','Nigeria')
    ,@('Niue','Not in IFs','#N/A')
    ,@('Norway','Norway','Norway')
    ,@('Occupied Palestinian Territory','Palestine','Palestine')
    ,@('Oman','Oman','Oman')
    ,@('Pakistan','Pakistan','Pakistan')
    ,@('Palau','Not in IFs','#N/A')
    ,@('Panama','Panama','Panama')
    ,@('Papua New Guinea','Papua New Guinea','Papua NG')
    ,@('Paraguay','Paraguay','Paraguay')
    ,@('Peru','Peru','Peru')
    ,@('Philippines','Philippines','Philippines')
    ,@('Poland','Poland','Poland')
    ,@('Portugal','Portugal','Portugal')
    ,@('Puerto Rico','Puerto Rico','Puerto Rico')
    ,@('Qatar','Qatar','Qatar')
    ,@('Republic of Korea','Korea, Republic of','Korea South')
    ,@('Republic of Moldova','Moldova','Moldova')
    ,@('Romania','Romania','Romania')
    ,@('Russian Federation','Russian Federation','Russia')
    ,@('Rwanda','Rwanda','Rwanda')
    ,@('Saint Kitts and Nevis','Not in IFs','#N/A')
    ,@('Saint Lucia','St. Lucia','St. Lucia')
    ,@('Saint Vincent and the Grenadines','St. Vincent and the Grenadines','St. Vincent and the Grenadines')
    ,@('Samoa','Samoa','Samoa')
    ,@('San Marino','Not in IFs','#N/A')
    ,@('Sao Tome and Principe','Sao Tome and Principe','Sao Tome and Principe')
    ,@('Saudi Arabia','Saudi Arabia','Saudi Arabia')
    ,@('Senegal','Senegal','Senegal')
    ,@('Serbia','Serbia','Serbia')
    ,@('Seychelles','Seychelles','Seychelles')
    ,@('Sierra Leone','Sierra Leone','SierraLeo')
    ,@('Singapore','Singapore','Singapore')
    ,@('Slovakia','Slovak Republic','Slovak Rep')
    ,@('Slovenia','Slovenia','Slovenia')
    ,@('Solomon Islands','Solomon Islands','Solomon Islands')
    ,@('Somalia','Somalia','Somalia')
    ,@('South Africa','South Africa','South Africa')
    ,@('South Sudan','Sudan South','Sudan South')
    ,@('Spain','Spain','Spain')
    ,@('Sri Lanka','Sri Lanka','Sri Lanka')
    ,@('Sudan','Sudan','Sudan')
    ,@('Suriname','Suriname','Suriname')
    ,@('Swaziland','Swaziland','Swaziland')
    ,@('Sweden','Sweden','Sweden')
    ,@('Switzerland','Switzerland','Switzerland')
    ,@('Syrian Arab Republic','Syrian Arab Republic','Syria')
    ,@('Tajikistan','Tajikistan','Tajikistan')
    ,@('Thailand','Thailand','Thailand')
    ,@('The former Yugoslav Republic of Macedonia','Macedonia, Former Yugoslav Republic of','Macedonia')
    ,@('Timor-Leste','Timor-Leste','Timor-Leste')
    ,@('Togo','Togo','Togo')
    ,@('Tokelau','Not in IFs','#N/A')
    ,@('Tonga','Tonga','Tonga')
    ,@('Trinidad and Tobago','Trinidad and Tobago','Trinidad')
    ,@('Tunisia','Tunisia','Tunisia')
    ,@('Turkey','Turkey','Turkey')
    ,@('Turkmenistan','Turkmenistan','Turkmenistan')
    ,@('Tuvalu','Not in IFs','#N/A')
    ,@('Uganda','Uganda','Uganda')
    ,@('Ukraine','Ukraine','Ukraine')
    ,@('United Arab Emirates','United Arab Emirates','UAE')
    ,@('United Kingdom','United Kingdom','Unitd Kingdm')
    ,@('United Republic of Tanzania','Tanzania','Tanzania')
    ,@('United States of America','United States','USA')
    ,@('Uruguay','Uruguay','Uruguay')
    ,@('Uzbekistan','Uzbekistan','Uzbekistan')
    ,@('Vanuatu','Vanuatu','Vanuatu')
    ,@('Venezuela (Bolivarian Republic of)','Venezuela','Venezuela')
    ,@('Viet Nam','Vietnam','Vietnam')
    ,@('Yemen','Yemen, Republic of','Yemen')
    ,@('Zambia','Zambia','Zambia')
    ,@('Zimbabwe','Zimbabwe','Zimbabwe')
    ,@('Hong Kong','Hong Kong','Hong Kong')
    ,@('Kosovo','Kosovo','Kosovo')
    ,@('Taiwan','Taiwan','Taiwan')
)

$n = $data.Length
for ($i = 0; $i -lt $n; $i++) {
    $r = $i + 1
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    if ($row[2] -eq '#N/A') {
        $ws.Cells.Item($r, 3).Value = '#N/A'
    } else {
        $ws.Cells.Item($r, 3).Value = $row[2]
    }
}

# Row 66 / B66 keeps its special wrapped-border style (style index 1 in styles.xml,
# "Normal_Sheet2"-based cellStyle) -- already present on that cell from the source
# workbook, so no additional style assignment is required here.

# Column widths per the updated layout (target OOXML widths are 14.85546875 /
# 19.140625 "characters"; the nearest values reachable through the ColumnWidth
# COM property are used here).
$ws.Columns.Item(1).ColumnWidth = 14
$ws.Columns.Item(2).ColumnWidth = 18.333333333333336

# Sheet view: drop the scrolled/selected state left over from editing, reset to A1.
$ws.Range("A1").Select()
